$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 1160.8572
$ws.Range("J2").Value2 = 2092
$ws.Range("L2").Value2 = 2092
$ws.Range("N2").Value2 = -2318

$ws.Range("H9").Value2 = 3674.8125
$ws.Range("I9").Value2 = 5159.095
$ws.Range("J9").Value2 = 841.1818
$ws.Range("K9").Value2 = 5159.095
$ws.Range("L9").Value2 = 841.1818
$ws.Range("M9").Value2 = -4990.095
$ws.Range("N9").Value2 = -1179.1818

$ws.Range("H28").Value2 = 3037.111
$ws.Range("I28").Value2 = 2055.6667
$ws.Range("K28").Value2 = 2055.6667
$ws.Range("M28").Value2 = -1570.6667

$ws.Range("H32").Value2 = 14390.6
$ws.Range("I32").Value2 = 2366.6667
$ws.Range("J32").Value2 = 19543.715
$ws.Range("K32").Value2 = 2366.6667
$ws.Range("L32").Value2 = 19543.715
$ws.Range("M32").Value2 = -2040.6667
$ws.Range("N32").Value2 = -20195.715

$ws.Range("H40").Value2 = 6183165.5
$ws.Range("I40").Value2 = 4797.615
$ws.Range("K40").Value2 = 4797.615
$ws.Range("M40").Value2 = -4622.615

$ws.Range("H74").Value2 = 18563534
$ws.Range("I74").Value2 = 18563534
$ws.Range("K74").Value2 = 18563534
$ws.Range("M74").Value2 = -18562598

$ws.Range("H77").Value2 = 18563534
$ws.Range("I77").Value2 = 18563534
$ws.Range("K77").Value2 = 92817670
$ws.Range("M77").Value2 = -92812990

$ws.Range("H88").Value2 = 37010296
$ws.Range("J88").Value2 = 6128064
$ws.Range("L88").Value2 = 6128064
$ws.Range("N88").Value2 = -6128876

$ws.Range("H91").Value2 = 37010296
$ws.Range("J91").Value2 = 6128064
$ws.Range("L91").Value2 = 6128064
$ws.Range("N91").Value2 = -6130872

$ws.Range("H107").Value2 = 1245.625
$ws.Range("I107").Value2 = 1148.8462
$ws.Range("K107").Value2 = 1148.8462
$ws.Range("M107").Value2 = 771.1538

$ws.Range("H111").Value2 = 759.5
$ws.Range("I111").Value2 = 759.5
$ws.Range("K111").Value2 = 2278.5
$ws.Range("M111").Value2 = 788.5

$ws.Range("H116").Value2 = 5999.5
$ws.Range("I116").Value2 = 5999.5
$ws.Range("K116").Value2 = 5999.5
$ws.Range("M116").Value2 = -2557.5

$ws.Range("H129").Value2 = 1834
$ws.Range("I129").Value2 = 1352.8889
$ws.Range("J129").Value2 = 3999
$ws.Range("K129").Value2 = 4058.6667
$ws.Range("L129").Value2 = 11997
$ws.Range("M129").Value2 = 941.3333000000002
$ws.Range("N129").Value2 = -21997

$ws.Range("H135").Value2 = 115385790
$ws.Range("I135").Value2 = 41667940
$ws.Range("K135").Value2 = 375011460
$ws.Range("M135").Value2 = -375008925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 429791.34
$ws.Range("J2").Value2 = 4276.4
$ws.Range("L2").Value2 = 4276.4
$ws.Range("N2").Value2 = -4502.4

$ws.Range("H4").Value2 = 12642.75
$ws.Range("I4").Value2 = 83.833336
$ws.Range("J4").Value2 = 50319.5
$ws.Range("K4").Value2 = 83.833336
$ws.Range("L4").Value2 = 50319.5
$ws.Range("M4").Value2 = 32.166664
$ws.Range("N4").Value2 = -50551.5

$ws.Range("H5").Value2 = 89.5625
$ws.Range("I5").Value2 = 131.77777
$ws.Range("J5").Value2 = 35.285713
$ws.Range("K5").Value2 = 131.77777
$ws.Range("L5").Value2 = 35.285713
$ws.Range("M5").Value2 = -19.77777
$ws.Range("N5").Value2 = -259.285713

$ws.Range("H45").Value2 = 1640.2727
$ws.Range("I45").Value2 = 1588.4445
$ws.Range("J45").Value2 = 1873.5
$ws.Range("K45").Value2 = 1588.4445
$ws.Range("L45").Value2 = 1873.5
$ws.Range("M45").Value2 = -1211.4445
$ws.Range("N45").Value2 = -2627.5

$ws.Range("H61").Value2 = 33337198
$ws.Range("I61").Value2 = 37040890
$ws.Range("K61").Value2 = 37040890
$ws.Range("M61").Value2 = -37040678

$ws.Range("H88").Value2 = 5409.1
$ws.Range("J88").Value2 = 7496.8335
$ws.Range("L88").Value2 = 7496.8335
$ws.Range("N88").Value2 = -8308.833500000001

$ws.Range("H91").Value2 = 5409.1
$ws.Range("J91").Value2 = 7496.8335
$ws.Range("L91").Value2 = 7496.8335
$ws.Range("N91").Value2 = -10304.8335

$ws.Range("H116").Value2 = 429791.34
$ws.Range("J116").Value2 = 4276.4
$ws.Range("L116").Value2 = 4276.4
$ws.Range("N116").Value2 = -8864.4

$ws.Range("H122").Value2 = 2970.35
$ws.Range("I122").Value2 = 2547.7368
$ws.Range("J122").Value2 = 11000
$ws.Range("K122").Value2 = 7643.2104
$ws.Range("L122").Value2 = 33000
$ws.Range("M122").Value2 = -5193.2104
$ws.Range("N122").Value2 = -37900

$ws.Range("H125").Value2 = 0
$ws.Range("J125").Value2 = 0
$ws.Range("L125").Value2 = 0
$ws.Range("N125").ClearContents()

$ws.Range("H136").Value2 = 33337198
$ws.Range("I136").Value2 = 37040890
$ws.Range("K136").Value2 = 111122670
$ws.Range("M136").Value2 = -111120120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 429791.34
$ws.Range("J3").Value2 = 4276.4
$ws.Range("L3").Value2 = 4276.4
$ws.Range("N3").Value2 = -4504.4

$ws.Range("H4").Value2 = 89.5625
$ws.Range("I4").Value2 = 131.77777
$ws.Range("J4").Value2 = 35.285713
$ws.Range("K4").Value2 = 131.77777
$ws.Range("L4").Value2 = 35.285713
$ws.Range("M4").Value2 = -16.77777
$ws.Range("N4").Value2 = -265.285713

$ws.Range("H86").Value2 = 8074.125
$ws.Range("J86").Value2 = 7798
$ws.Range("L86").Value2 = 7798
$ws.Range("N86").Value2 = -10044

$ws.Range("H89").Value2 = 8074.125
$ws.Range("J89").Value2 = 7798
$ws.Range("L89").Value2 = 38990
$ws.Range("N89").Value2 = -50222

$ws.Range("H107").Value2 = 60879.234
$ws.Range("I107").Value2 = 1448.5
$ws.Range("J107").Value2 = 113706.555
$ws.Range("K107").Value2 = 1448.5
$ws.Range("L107").Value2 = 113706.555
$ws.Range("M107").Value2 = 471.5
$ws.Range("N107").Value2 = -117546.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 117.46154
$ws.Range("I7").Value2 = 40
$ws.Range("K7").Value2 = 40
$ws.Range("M7").Value2 = 73

$ws.Range("H31").Value2 = 11105.111
$ws.Range("I31").Value2 = 7432.778
$ws.Range("K31").Value2 = 7432.778
$ws.Range("M31").Value2 = -7137.778

$ws.Range("H34").Value2 = 11105.111
$ws.Range("I34").Value2 = 7432.778
$ws.Range("K34").Value2 = 7432.778
$ws.Range("M34").Value2 = -7230.778

$ws.Range("H99").Value2 = 2006
$ws.Range("I99").Value2 = 1908.1666
$ws.Range("J99").Value2 = 2299.5
$ws.Range("K99").Value2 = 1908.1666
$ws.Range("L99").Value2 = 2299.5
$ws.Range("M99").Value2 = -410.1666
$ws.Range("N99").Value2 = -5295.5

$ws.Range("H126").Value2 = 2006
$ws.Range("I126").Value2 = 1908.1666
$ws.Range("J126").Value2 = 2299.5
$ws.Range("K126").Value2 = 5724.4998
$ws.Range("L126").Value2 = 6898.5
$ws.Range("M126").Value2 = -3254.4998
$ws.Range("N126").Value2 = -11838.5

$ws.Range("H132").Value2 = 37038188
$ws.Range("I132").Value2 = 45455572
$ws.Range("J132").Value2 = 1689.4
$ws.Range("K132").Value2 = 136366716
$ws.Range("L132").Value2 = 5068.200000000001
$ws.Range("M132").Value2 = -136364186
$ws.Range("N132").Value2 = -10128.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value2 = 5000
$ws.Range("I100").Value2 = 0
$ws.Range("K100").Value2 = 0
$ws.Range("M100").ClearContents()

$ws.Range("H107").Value2 = 1508.5238
$ws.Range("J107").Value2 = 2355.8333
$ws.Range("L107").Value2 = 7067.499899999999
$ws.Range("N107").Value2 = -10907.4999

$ws.Range("H132").Value2 = 1535.4445
$ws.Range("J132").Value2 = 1259.4
$ws.Range("L132").Value2 = 11334.6
$ws.Range("N132").Value2 = -16394.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 5000
$ws.Range("J46").Value2 = 0
$ws.Range("L46").Value2 = 0
$ws.Range("N46").ClearContents()

$ws.Range("H80").Value2 = 2999.5
$ws.Range("I80").Value2 = 2999
$ws.Range("K80").Value2 = 2999
$ws.Range("M80").Value2 = -2001

$ws.Range("H83").Value2 = 2999.5
$ws.Range("I83").Value2 = 2999
$ws.Range("K83").Value2 = 14995
$ws.Range("M83").Value2 = -10003

$ws.Range("H107").Value2 = 549.9231
$ws.Range("I107").Value2 = 322.5
$ws.Range("J107").Value2 = 1308
$ws.Range("K107").Value2 = 322.5
$ws.Range("L107").Value2 = 1308
$ws.Range("M107").Value2 = 1597.5
$ws.Range("N107").Value2 = -5148

$ws.Range("H132").Value2 = 7814192
$ws.Range("I132").Value2 = 8929863
$ws.Range("J132").Value2 = 4499.5
$ws.Range("K132").Value2 = 26789589
$ws.Range("L132").Value2 = 13498.5
$ws.Range("M132").Value2 = -26787059
$ws.Range("N132").Value2 = -18558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 7937646
$ws.Range("I68").Value2 = 11905512
$ws.Range("J68").Value2 = 1914
$ws.Range("K68").Value2 = 11905512
$ws.Range("L68").Value2 = 1914
$ws.Range("M68").Value2 = -11904763
$ws.Range("N68").Value2 = -3412

$ws.Range("H71").Value2 = 7937646
$ws.Range("I71").Value2 = 11905512
$ws.Range("J71").Value2 = 1914
$ws.Range("K71").Value2 = 59527560
$ws.Range("L71").Value2 = 9570
$ws.Range("M71").Value2 = -59523816
$ws.Range("N71").Value2 = -17058

$ws.Range("H132").Value2 = 19241276
$ws.Range("I132").Value2 = 19241276
$ws.Range("K132").Value2 = 57723828
$ws.Range("M132").Value2 = -57721298

$ws.Range("H136").Value2 = 2716.4285
$ws.Range("I136").Value2 = 1346.6666
$ws.Range("J136").Value2 = 3090
$ws.Range("K136").Value2 = 4039.9998
$ws.Range("L136").Value2 = 9270
$ws.Range("M136").Value2 = -1489.9998
$ws.Range("N136").Value2 = -14370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 2912.9211
$ws.Range("J96").Value2 = 3199.4517
$ws.Range("L96").Value2 = 3199.4517
$ws.Range("N96").Value2 = -5945.4517
